$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion-of-the-day message text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.25 = 29246.38 pesos`n✅ 29246.38 pesos = 7.24 = 942.32 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Sheet "tasas": update the rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 138
$ws2.Range("O10").Value = 4036
$ws2.Range("N12").Value = 4041
$ws2.Range("O12").Value = 130.201
